# FR24 bug report workbook update:
#  - Fill in "Date Fixed" (B column) for the bugs that are now resolved
#  - Normalize the Severity/Priority/Reported By cell formatting
#    (drop the stray "apply fill" flag those cells had picked up)
#  - Restore the sheet's cursor/selection to B2

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ZK_Bug_Report")

# "Date Fixed" rows for each bug block now have a fix date recorded
# (2023-12-11). B3 already carries the exact date style (yyyy-mm-dd,
# left/center, wrapped) we need, so copy its formatting across and then
# overwrite just the value.
$dateStyleSource = $ws.Range("B3")
$fixedRows = 16, 34, 52, 70, 88, 106, 124, 142
foreach ($r in $fixedRows) {
    $cell = $ws.Cells.Item($r, 2)
    $dateStyleSource.Copy($cell)
    $cell.Value = 45271
}

# Severity / Priority / Reported By cells for bugs 5, 6, 7 and 8 had an
# extraneous fill-applied style; clear it so they match the rest of the
# template (no fill).
$fillFixRows = 77, 78, 79, 95, 96, 97, 113, 114, 115, 131, 132, 133
foreach ($r in $fillFixRows) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Interior.Pattern = -4142
}

$ws.Activate() | Out-Null
$ws.Range("B2").Select() | Out-Null
